$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule_date")

# Module 2's label is updated from "2: Coding fundamental" to
# "2: Coding fundamentals" for the three rows that use it (first
# addition of the "033 R Markdown" module content).
$ws.Range("C4").Value = "2: Coding fundamentals"
$ws.Range("C5").Value = "2: Coding fundamentals"
$ws.Range("C6").Value = "2: Coding fundamentals"

# Move the saved selection to D7, as recorded in the sheet view.
$ws.Activate()
$ws.Range("D7").Select()
